# Fix a typo in the target-address column and move the selection,
# matching the authored diff:
#   - C6: "192.168.0.12/24" -> "192.168.12.0/24"
#   - selection moves from E5 to C7 (and the saved topLeftCell scroll anchor is cleared)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = "192.168.12.0/24"

$ws.Range("C7").Select()
